# datannur __metaVariable__.xlsx
# "fixed : tag level on tab tag page dataset, added : dataset delivery_format"
#
# The "dataset" entity is missing a "delivery_format" variable row (the
# "folder" entity already documents one). Insert it right after the
# "localisation" / dataset row (new row 20), which keeps the Tableau3
# table sorted by the "dataset" (entity) column, and shifts every
# following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 20 - everything from the old row 20
# ("start_date" / dataset) down to the old last row (109 -> 110 -> 111)
# shifts down by one.
$ws.Rows("20:20").Insert()

# Fill in the new "delivery_format" row for the "dataset" entity.
$ws.Range("A20").Value = "delivery_format"
$ws.Range("B20").Value = "dataset"
$ws.Range("C20").Value = "Format du dataset livrées (CSV, XML, ...)"

# Grow the Tableau3 table (was A1:C110) so it covers the new row and the
# row that got pushed out to the new last position (A1:C111).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C111"))

# Leave the selection where the author left it after the edit.
$ws.Range("C21").Select()
